$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1658.8148
$ws.Range("J17").Value = 1856.5714
$ws.Range("L17").Value = 5569.7142
$ws.Range("N17").Value = -5905.7142
$ws.Range("H40").Value = 6184
$ws.Range("I40").Value = 5141.25
$ws.Range("K40").Value = 5141.25
$ws.Range("M40").Value = -4966.25
$ws.Range("H80").Value = 569.55554
$ws.Range("I80").Value = 505.16666
$ws.Range("J80").Value = 698.3333
$ws.Range("K80").Value = 1515.49998
$ws.Range("L80").Value = 2094.9999
$ws.Range("M80").Value = -517.4999800000001
$ws.Range("N80").Value = -4090.9999
$ws.Range("H83").Value = 569.55554
$ws.Range("I83").Value = 505.16666
$ws.Range("J83").Value = 698.3333
$ws.Range("K83").Value = 4546.49994
$ws.Range("L83").Value = 6284.9997
$ws.Range("M83").Value = 445.5000600000003
$ws.Range("N83").Value = -16268.9997
$ws.Range("H94").Value = 3458.3333
$ws.Range("I94").Value = 3454.5454
$ws.Range("K94").Value = 3454.5454
$ws.Range("M94").Value = -3003.5454
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = $null
$ws.Range("N123").Value = 0
$ws.Range("H131").Value = 1211.8
$ws.Range("I131").Value = 1211.8
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 3635.4
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = $null
$ws.Range("N131").Value = 1404.6
$ws.Range("H132").Value = 15749.765
$ws.Range("I132").Value = 15382.214
$ws.Range("K132").Value = 46146.642
$ws.Range("M132").Value = -43616.642
$ws.Range("H137").Value = 3810.1428
$ws.Range("I137").Value = 859.3333
$ws.Range("K137").Value = 2577.9999
$ws.Range("M137").Value = -27.9998999999998
$ws.Range("H138").Value = 2484.2856
$ws.Range("I138").Value = 1143.6364
$ws.Range("J138").Value = 7400
$ws.Range("K138").Value = 3430.9092
$ws.Range("L138").Value = 22200
$ws.Range("M138").Value = 1709.0908
$ws.Range("N138").Value = -32480
$ws.Range("H141").Value = 963.25
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = $null

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1500
$ws.Range("I2").Value = 1500
$ws.Range("K2").Value = 1500
$ws.Range("M2").Value = -1387
$ws.Range("H32").Value = 3634.5625
$ws.Range("I32").Value = 604.4815
$ws.Range("K32").Value = 604.4815
$ws.Range("M32").Value = -317.4815
$ws.Range("H74").Value = 3418.0454
$ws.Range("I74").Value = 3105.1
$ws.Range("K74").Value = 3105.1
$ws.Range("M74").Value = -2231.1
$ws.Range("H77").Value = 3418.0454
$ws.Range("I77").Value = 3105.1
$ws.Range("K77").Value = 15525.5
$ws.Range("M77").Value = -11157.5
$ws.Range("H88").Value = 904.5
$ws.Range("I88").Value = 397.33334
$ws.Range("J88").Value = 1665.25
$ws.Range("K88").Value = 397.33334
$ws.Range("L88").Value = 1665.25
$ws.Range("M88").Value = 8.666659999999979
$ws.Range("N88").Value = -2477.25
$ws.Range("H91").Value = 904.5
$ws.Range("I91").Value = 397.33334
$ws.Range("J91").Value = 1665.25
$ws.Range("K91").Value = 397.33334
$ws.Range("L91").Value = 1665.25
$ws.Range("M91").Value = 1006.66666
$ws.Range("N91").Value = -4473.25
$ws.Range("H116").Value = 1500
$ws.Range("I116").Value = 1500
$ws.Range("K116").Value = 1500
$ws.Range("M116").Value = 794

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1500
$ws.Range("I3").Value = 1500
$ws.Range("K3").Value = 1500
$ws.Range("M3").Value = -1386
$ws.Range("H105").Value = 1964
$ws.Range("I105").Value = 1587
$ws.Range("K105").Value = 1587
$ws.Range("M105").Value = 160

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 903.1539
$ws.Range("I22").Value = 887
$ws.Range("K22").Value = 887
$ws.Range("M22").Value = -537
$ws.Range("H31").Value = 4514
$ws.Range("J31").Value = 5168.4375
$ws.Range("L31").Value = 5168.4375
$ws.Range("N31").Value = -5758.4375
$ws.Range("H34").Value = 4514
$ws.Range("J34").Value = 5168.4375
$ws.Range("L34").Value = 5168.4375
$ws.Range("N34").Value = -5572.4375
$ws.Range("H51").Value = 39129.8
$ws.Range("J51").Value = 45588.332
$ws.Range("L51").Value = 45588.332
$ws.Range("N51").Value = -47060.332
$ws.Range("H58").Value = 3037.9546
$ws.Range("I58").Value = 1769.7222
$ws.Range("K58").Value = 1769.7222
$ws.Range("M58").Value = -1566.7222
$ws.Range("H60").Value = 22849.9
$ws.Range("H61").Value = 39129.8
$ws.Range("J61").Value = 45588.332
$ws.Range("L61").Value = 45588.332
$ws.Range("N61").Value = -46284.332
$ws.Range("H122").Value = 583
$ws.Range("I122").Value = 525
$ws.Range("K122").Value = 1575
$ws.Range("M122").Value = 875
$ws.Range("H132").Value = 1972.62
$ws.Range("I132").Value = 2009.4651
$ws.Range("K132").Value = 6028.3953
$ws.Range("M132").Value = -3498.3953
$ws.Range("H134").Value = 2049.1428
$ws.Range("I134").Value = 1204.4
$ws.Range("J134").Value = 4161
$ws.Range("K134").Value = 3613.2
$ws.Range("L134").Value = 12483
$ws.Range("M134").Value = -1078.2
$ws.Range("N134").Value = -17553
$ws.Range("H136").Value = 3037.9546
$ws.Range("I136").Value = 1769.7222
$ws.Range("K136").Value = 5309.1666
$ws.Range("M136").Value = -2759.1666

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2166.6667
$ws.Range("J80").Value = 2250
$ws.Range("L80").Value = 2250
$ws.Range("N80").Value = -4246
$ws.Range("H83").Value = 2166.6667
$ws.Range("J83").Value = 2250
$ws.Range("L83").Value = 11250
$ws.Range("N83").Value = -21234
$ws.Range("H132").Value = 32034.914
$ws.Range("I132").Value = 36640.934
$ws.Range("J132").Value = 4398.8
$ws.Range("K132").Value = 109922.802
$ws.Range("L132").Value = 13196.4
$ws.Range("M132").Value = -107392.802
$ws.Range("N132").Value = -18256.4

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = $null
$ws.Range("M16").Value = $null
$ws.Range("N16").Value = 0
$ws.Range("H22").Value = 1380.4375
$ws.Range("I22").Value = 778.8
$ws.Range("J22").Value = 1653.909
$ws.Range("K22").Value = 778.8
$ws.Range("L22").Value = 1653.909
$ws.Range("M22").Value = -483.8
$ws.Range("N22").Value = -2243.909
$ws.Range("H27").Value = 1380.4375
$ws.Range("I27").Value = 778.8
$ws.Range("J27").Value = 1653.909
$ws.Range("K27").Value = 778.8
$ws.Range("L27").Value = 1653.909
$ws.Range("M27").Value = -671.8
$ws.Range("N27").Value = -1867.909
$ws.Range("H46").Value = 1871.5
$ws.Range("J46").Value = 1001.75
$ws.Range("L46").Value = 1001.75
$ws.Range("N46").Value = -1377.75
$ws.Range("H132").Value = 4120.4165
$ws.Range("I132").Value = 2160
$ws.Range("K132").Value = 6480
$ws.Range("M132").Value = -3950
$ws.Range("H136").Value = 4260
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").Value = $null

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 60000
$ws.Range("J64").Value = 60000
$ws.Range("L64").Value = 60000
$ws.Range("N64").Value = -60496
$ws.Range("H67").Value = 60000
$ws.Range("J67").Value = 60000
$ws.Range("L67").Value = 60000
$ws.Range("N67").Value = -61716
$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800
$ws.Range("H132").Value = 2658.7273
$ws.Range("I132").Value = 2093.25
$ws.Range("J132").Value = 4166.6665
$ws.Range("K132").Value = 6279.75
$ws.Range("L132").Value = 12499.9995
$ws.Range("M132").Value = -3749.75
$ws.Range("N132").Value = -17559.9995
$ws.Range("H136").Value = 2440.7778
$ws.Range("I136").Value = 1935.6
$ws.Range("K136").Value = 5806.799999999999
$ws.Range("M136").Value = -3256.799999999999
